$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 28.89432766666667
$ws.Range("H2").Value = 86.68298300000001
$ws.Range("I2").Value = 0.9344268072004271
$ws.Range("J2").Value = 0.934426807200427
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.642196
$ws.Range("N2").Value = 7.926588000000001
$ws.Range("O2").Value = 0.26568831615543
$ws.Range("P2").Value = 0.26568831615543
$ws.Range("Q2").Value = 76.344476983556
$ws.Range("R2").Value = 687.1002928520041
$ws.Range("S2").Value = 0.2482662849755761
$ws.Range("T2").Value = 0.2482662849755761

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 28.89432766666667
$ws.Range("H3").Value = 86.68298300000001
$ws.Range("I3").Value = 0.9344268072004271
$ws.Range("J3").Value = 0.934426807200427
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.666004666666667
$ws.Range("N3").Value = 13.998014
$ws.Range("O3").Value = 0.469194156323015
$ws.Range("P3").Value = 0.4691941563230151
$ws.Range("Q3").Value = 134.8210677328624
$ws.Range("R3").Value = 1213.389609595762
$ws.Range("S3").Value = 0.438427597450013
$ws.Range("T3").Value = 0.438427597450013

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 28.89432766666667
$ws.Range("H4").Value = 86.68298300000001
$ws.Range("I4").Value = 0.9344268072004271
$ws.Range("J4").Value = 0.934426807200427
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.114591666666667
$ws.Range("N4").Value = 3.343775
$ws.Range("O4").Value = 0.1120787341732184
$ws.Range("P4").Value = 0.1120787341732184
$ws.Range("Q4").Value = 32.20537683120278
$ws.Range("R4").Value = 289.848391480825
$ws.Range("S4").Value = 0.1047293737285459
$ws.Range("T4").Value = 0.1047293737285459

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 28.89432766666667
$ws.Range("H5").Value = 86.68298300000001
$ws.Range("I5").Value = 0.9344268072004271
$ws.Range("J5").Value = 0.934426807200427
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.521928
$ws.Range("N5").Value = 4.565784000000001
$ws.Range("O5").Value = 0.1530387933483365
$ws.Range("P5").Value = 0.1530387933483365
$ws.Range("Q5").Value = 43.97508631707468
$ws.Range("R5").Value = 395.7757768536721
$ws.Range("S5").Value = 0.1430035510462921
$ws.Range("T5").Value = 0.143003551046292

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.1893213333333333
$ws.Range("H6").Value = 0.567964
$ws.Range("I6").Value = 0.006122548725910637
$ws.Range("J6").Value = 0.006122548725910637
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.642196
$ws.Range("N6").Value = 7.926588000000001
$ws.Range("O6").Value = 0.26568831615543
$ws.Range("P6").Value = 0.26568831615543
$ws.Range("Q6").Value = 0.500224069648
$ws.Range("R6").Value = 4.502016626832001
$ws.Range("S6").Value = 0.00162668966156677
$ws.Range("T6").Value = 0.001626689661566771

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.1893213333333333
$ws.Range("H7").Value = 0.567964
$ws.Range("I7").Value = 0.006122548725910637
$ws.Range("J7").Value = 0.006122548725910637
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.666004666666667
$ws.Range("N7").Value = 13.998014
$ws.Range("O7").Value = 0.469194156323015
$ws.Range("P7").Value = 0.4691941563230151
$ws.Range("Q7").Value = 0.883374224832889
$ws.Range("R7").Value = 7.950368023496001
$ws.Range("S7").Value = 0.002872664084000192
$ws.Range("T7").Value = 0.002872664084000192

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.1893213333333333
$ws.Range("H8").Value = 0.567964
$ws.Range("I8").Value = 0.006122548725910637
$ws.Range("J8").Value = 0.006122548725910637
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.114591666666667
$ws.Range("N8").Value = 3.343775
$ws.Range("O8").Value = 0.1120787341732184
$ws.Range("P8").Value = 0.1120787341732184
$ws.Range("Q8").Value = 0.2110159804555556
$ws.Range("R8").Value = 1.8991438241
$ws.Range("S8").Value = 0.0006862075111139152
$ws.Range("T8").Value = 0.0006862075111139152

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.1893213333333333
$ws.Range("H9").Value = 0.567964
$ws.Range("I9").Value = 0.006122548725910637
$ws.Range("J9").Value = 0.006122548725910637
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.521928
$ws.Range("N9").Value = 4.565784000000001
$ws.Range("O9").Value = 0.1530387933483365
$ws.Range("P9").Value = 0.1530387933483365
$ws.Range("Q9").Value = 0.2881334381973334
$ws.Range("R9").Value = 2.593200943776
$ws.Range("S9").Value = 0.0009369874692297589
$ws.Range("T9").Value = 0.0009369874692297589

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.339639
$ws.Range("H10").Value = 4.018917
$ws.Range("I10").Value = 0.04332319505794487
$ws.Range("J10").Value = 0.04332319505794486
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.642196
$ws.Range("N10").Value = 7.926588000000001
$ws.Range("O10").Value = 0.26568831615543
$ws.Range("P10").Value = 0.26568831615543
$ws.Range("Q10").Value = 3.539588807244
$ws.Range("R10").Value = 31.856299265196
$ws.Range("S10").Value = 0.01151046674541862
$ws.Range("T10").Value = 0.01151046674541862

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.339639
$ws.Range("H11").Value = 4.018917
$ws.Range("I11").Value = 0.04332319505794487
$ws.Range("J11").Value = 0.04332319505794486
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.666004666666667
$ws.Range("N11").Value = 13.998014
$ws.Range("O11").Value = 0.469194156323015
$ws.Range("P11").Value = 0.4691941563230151
$ws.Range("Q11").Value = 6.250761825648667
$ws.Range("R11").Value = 56.25685643083801
$ws.Range("S11").Value = 0.02032698995442985
$ws.Range("T11").Value = 0.02032698995442985

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.339639
$ws.Range("H12").Value = 4.018917
$ws.Range("I12").Value = 0.04332319505794487
$ws.Range("J12").Value = 0.04332319505794486
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.114591666666667
$ws.Range("N12").Value = 3.343775
$ws.Range("O12").Value = 0.1120787341732184
$ws.Range("P12").Value = 0.1120787341732184
$ws.Range("Q12").Value = 1.493150465741667
$ws.Range("R12").Value = 13.438354191675
$ws.Range("S12").Value = 0.004855608862433892
$ws.Range("T12").Value = 0.004855608862433891

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.339639
$ws.Range("H13").Value = 4.018917
$ws.Range("I13").Value = 0.04332319505794487
$ws.Range("J13").Value = 0.04332319505794486
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.521928
$ws.Range("N13").Value = 4.565784000000001
$ws.Range("O13").Value = 0.1530387933483365
$ws.Range("P13").Value = 0.1530387933483365
$ws.Range("Q13").Value = 2.038834103992
$ws.Range("R13").Value = 18.349506935928
$ws.Range("S13").Value = 0.006630129495662497
$ws.Range("T13").Value = 0.006630129495662497

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4986926666666666
$ws.Range("H14").Value = 1.496078
$ws.Range("I14").Value = 0.01612744901571743
$ws.Range("J14").Value = 0.01612744901571742
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.642196
$ws.Range("N14").Value = 7.926588000000001
$ws.Range("O14").Value = 0.26568831615543
$ws.Range("P14").Value = 0.26568831615543
$ws.Range("Q14").Value = 1.317643769096
$ws.Range("R14").Value = 11.858793921864
$ws.Range("S14").Value = 0.00428487477286851
$ws.Range("T14").Value = 0.00428487477286851

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4986926666666666
$ws.Range("H15").Value = 1.496078
$ws.Range("I15").Value = 0.01612744901571743
$ws.Range("J15").Value = 0.01612744901571742
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 4.666004666666667
$ws.Range("N15").Value = 13.998014
$ws.Range("O15").Value = 0.469194156323015
$ws.Range("P15").Value = 0.4691941563230151
$ws.Range("Q15").Value = 2.326902309899111
$ws.Range("R15").Value = 20.942120789092
$ws.Range("S15").Value = 0.007566904834571977
$ws.Range("T15").Value = 0.007566904834571976

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4986926666666666
$ws.Range("H16").Value = 1.496078
$ws.Range("I16").Value = 0.01612744901571743
$ws.Range("J16").Value = 0.01612744901571742
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.114591666666667
$ws.Range("N16").Value = 3.343775
$ws.Range("O16").Value = 0.1120787341732184
$ws.Range("P16").Value = 0.1120787341732184
$ws.Range("Q16").Value = 0.5558386904944443
$ws.Range("R16").Value = 5.002548214449999
$ws.Range("S16").Value = 0.001807544071124726
$ws.Range("T16").Value = 0.001807544071124726

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.4986926666666666
$ws.Range("H17").Value = 1.496078
$ws.Range("I17").Value = 0.01612744901571743
$ws.Range("J17").Value = 0.01612744901571742
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.521928
$ws.Range("N17").Value = 4.565784000000001
$ws.Range("O17").Value = 0.1530387933483365
$ws.Range("P17").Value = 0.1530387933483365
$ws.Range("Q17").Value = 0.7589743327946666
$ws.Range("R17").Value = 6.830768995152
$ws.Range("S17").Value = 0.002468125337152212
$ws.Range("T17").Value = 0.002468125337152212
